# Update completion percentages for Jump (up/over) and Shoot feature rows
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("F8").Value = 0.5
$ws.Range("F9").Value = 0.5
$ws.Range("F11").Value = 0.75
$ws.Range("F12").Value = 0.75
